$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = -3
    "F4"  = 0
    "F7"  = -2
    "F12" = 3
    "F14" = 1
    "F15" = 3
    "F17" = 1
    "F19" = 1
    "F22" = 0
    "F29" = -3
    "F37" = 3
    "F44" = -7
    "F48" = 4
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
